$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-25: update Price (D) and Volume(1h) (E) only where changed
$ws.Range("D2").Value = "60.206.08"
$ws.Range("E2").Value = "  +4.85%  "
$ws.Range("D3").Value = "2.598.14"
$ws.Range("E3").Value = "  +7.62%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "505.93"
$ws.Range("E5").Value = "  +3.61%  "
$ws.Range("D6").Value = "156.29"
$ws.Range("E6").Value = "  +1.91%  "
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  -3.54%  "
$ws.Range("D9").Value = "2.631.57"
$ws.Range("E9").Value = "  +8.10%  "
$ws.Range("D10").Value = "6.47"
$ws.Range("E10").Value = "  +6.00%  "
$ws.Range("E11").Value = "  +4.46%  "
$ws.Range("D12").Value = "0.342"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "3.078.00"
$ws.Range("E14").Value = "  +8.39%  "
$ws.Range("D15").Value = "60.313.56"
$ws.Range("E15").Value = "  +5.16%  "
$ws.Range("D16").Value = "21.65"
$ws.Range("E16").Value = "  +5.35%  "
$ws.Range("E17").Value = "  +5.04%  "
$ws.Range("D18").Value = "2.624.79"
$ws.Range("E18").Value = "  +7.76%  "
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("D20").Value = "343.46"
$ws.Range("E20").Value = "  +6.09%  "
$ws.Range("D21").Value = "10.44"
$ws.Range("E21").Value = "  +4.51%  "
$ws.Range("E22").Value = "  +3.89%  "
$ws.Range("D23").Value = "0.995"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "59.99"
$ws.Range("E24").Value = "  +3.55%  "
$ws.Range("E25").Value = "  +4.96%  "

# Rows 26-51: new coin inserted at 26, shifting ranking down; row 52 (former WhiteBITCoin) dropped
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.741.23"
$ws.Range("E26").Value = "  +8.66%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.166"
$ws.Range("E27").Value = "  +3.47%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "0.992"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0857"
$ws.Range("E29").Value = "  +9.36%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "7.54"
$ws.Range("E30").Value = "  +3.66%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "19.46"
$ws.Range("E32").Value = "  +4.51%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "156.11"
$ws.Range("E33").Value = "  +3.44%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "1.57"
$ws.Range("E34").Value = "  +3.20%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "5.74"
$ws.Range("E35").Value = "  +8.37%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "4.03"
$ws.Range("E36").Value = "  +7.13%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.20"
$ws.Range("E37").Value = "  +4.61%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "307.12"
$ws.Range("E38").Value = "  +8.33%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "1.49"
$ws.Range("E39").Value = "  +8.36%  "
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "0.847"
$ws.Range("E40").Value = "  +3.87%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "3.77"
$ws.Range("E41").Value = "  +7.10%  "
$ws.Range("B42").Value = "SuiNetwork"
$ws.Range("C42").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D42").Value = "0.836"
$ws.Range("E42").Value = "  +28.62%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "35.45"
$ws.Range("E43").Value = "  +4.36%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "0.629"
$ws.Range("E44").Value = "  +5.22%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0572"
$ws.Range("E45").Value = "  +7.92%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.100"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Value = "0.992"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "19.86"
$ws.Range("E48").Value = "  +12.84%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "4.88"
$ws.Range("E49").Value = "  +7.79%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.049.19"
$ws.Range("E50").Value = "  +7.86%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "0.0235"
$ws.Range("E51").Value = "  +3.31%  "
